$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.853.02'
$ws.Range('E2').Value = '  -2.23%  '
$ws.Range('D3').Value = '1.816.09'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').Value = '1.007'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('D6').Value = '308.38'
$ws.Range('E6').Value = '  -2.07%  '
$ws.Range('D7').Value = '0.4618'
$ws.Range('E7').Value = '  -2.59%  '
$ws.Range('D8').Value = '0.3634'
$ws.Range('E8').Value = '  -1.59%  '
$ws.Range('D9').Value = '0.07218'
$ws.Range('E9').Value = '  -3.26%  '
$ws.Range('D10').Value = '0.8574'
$ws.Range('E10').Value = '  -3.21%  '
$ws.Range('E11').Value = '  -3.59%  '
$ws.Range('D12').Value = '0.07515'
$ws.Range('E12').Value = '  +2.47%  '
$ws.Range('D13').Value = '1.768.94'
$ws.Range('E13').Value = '  -9.44%  '
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('D15').Value = '6.523'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D16').Value = '91.72'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '0.000008561'
$ws.Range('E18').Value = '  -2.93%  '
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').Value = '26.993.43'
$ws.Range('E20').Value = '  -1.81%  '
$ws.Range('D21').Value = '14.40'
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('E22').Value = '  -3.45%  '
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('D24').Value = '2.071.85'
$ws.Range('E24').Value = '  -2.94%  '
$ws.Range('D25').Value = '151.13'
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('E26').Value = '  -2.99%  '
$ws.Range('D27').Value = '18.14'
$ws.Range('E27').Value = '  -2.79%  '
$ws.Range('D28').Value = '2.063'
$ws.Range('E28').Value = '  -4.00%  '
$ws.Range('D29').Value = '5.085'
$ws.Range('E29').Value = '  -3.17%  '
$ws.Range('D30').Value = '115.00'
$ws.Range('E30').Value = '  -2.48%  '
$ws.Range('D31').Value = '0.08854'
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('D33').Value = '4.406'
$ws.Range('E33').Value = '  -3.39%  '
$ws.Range('E34').Value = '  -4.46%  '
$ws.Range('D35').Value = '0.7161'
$ws.Range('E35').Value = '  -5.38%  '
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('D37').Value = '1.073'
$ws.Range('E37').Value = '  -2.87%  '
$ws.Range('D38').Value = '2.438'
$ws.Range('E38').Value = '  +1.46%  '
$ws.Range('D39').Value = '0.05232'
$ws.Range('E39').Value = '  -1.93%  '
$ws.Range('E40').Value = '  -2.32%  '
$ws.Range('D41').Value = '2.920'
$ws.Range('E41').Value = '  -2.46%  '
$ws.Range('E42').Value = '  -2.72%  '
$ws.Range('D43').Value = '0.5127'
$ws.Range('E43').Value = '  -3.93%  '
$ws.Range('D44').Value = '0.1619'
$ws.Range('E44').Value = '  -2.51%  '
$ws.Range('D45').Value = '8.171'
$ws.Range('E45').Value = '  -4.05%  '
$ws.Range('D46').Value = '0.4789'
$ws.Range('E46').Value = '  -2.63%  '
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('D48').Value = '102.88'
$ws.Range('E48').Value = '  -2.23%  '
$ws.Range('D49').Value = '10.03'
$ws.Range('E49').Value = '  -4.75%  '
$ws.Range('D50').Value = '1.615'
$ws.Range('D51').Value = '0.06193'
$ws.Range('E51').Value = '  -2.04%  '
